$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Delete now-removed trailing rows (23-27), bottom-up ---
$ws.Rows.Item(27).Delete() | Out-Null
$ws.Rows.Item(26).Delete() | Out-Null
$ws.Rows.Item(25).Delete() | Out-Null
$ws.Rows.Item(24).Delete() | Out-Null
$ws.Rows.Item(23).Delete() | Out-Null

# --- Rewrite rows 10-22 (A/B/C) to match the restructured content ---
# Row 10
$ws.Cells.Item(10, 1).Value = 'Objetivos:'
$ws.Cells.Item(10, 2).Value = '471420 - Carlos Antonio Reis Pereira Baptista'
$ws.Cells.Item(10, 3).Value = '471420 - Carlos Antonio Reis Pereira Baptista'
$ws.Rows.Item(10).RowHeight = 60

# Row 11
$ws.Cells.Item(11, 1).Value = 'Objectives:'
$ws.Cells.Item(11, 2).ClearContents() | Out-Null
$ws.Cells.Item(11, 3).ClearContents() | Out-Null
$ws.Rows.Item(11).RowHeight = 60

# Row 12
$ws.Cells.Item(12, 1).Value = 'Programa resumido:'
$ws.Cells.Item(12, 2).Value = '5840897 - Clodoaldo Saron'
$ws.Cells.Item(12, 3).Value = '5840897 - Clodoaldo Saron'
$ws.Rows.Item(12).RowHeight = 60

# Row 13
$ws.Cells.Item(13, 1).Value = 'Short syllabus:'
$ws.Cells.Item(13, 2).ClearContents() | Out-Null
$ws.Cells.Item(13, 3).ClearContents() | Out-Null
$ws.Rows.Item(13).RowHeight = 60

# Row 14
$ws.Cells.Item(14, 1).Value = 'Programa:'
$ws.Cells.Item(14, 2).Value = '1033242 - Fábio Herbst Florenzano'
$ws.Cells.Item(14, 3).Value = '1033242 - Fábio Herbst Florenzano'
$ws.Rows.Item(14).RowHeight = 120

# Row 15
$ws.Cells.Item(15, 1).Value = 'Syllabus:'
$ws.Cells.Item(15, 2).ClearContents() | Out-Null
$ws.Cells.Item(15, 3).ClearContents() | Out-Null
$ws.Rows.Item(15).RowHeight = 120

# Row 16
$ws.Cells.Item(16, 1).Value = 'Avaliação:'
$ws.Cells.Item(16, 2).ClearContents() | Out-Null
$ws.Cells.Item(16, 3).ClearContents() | Out-Null

# Row 17
$ws.Cells.Item(17, 1).Value = 'Método:'
$ws.Cells.Item(17, 2).Value = '5840793 - Sérgio Schneider'
$ws.Cells.Item(17, 3).Value = '5840793 - Sérgio Schneider'
$ws.Rows.Item(17).RowHeight = 60

# Row 18
$ws.Cells.Item(18, 1).Value = 'Critério:'
$ws.Cells.Item(18, 2).Value = 'Experimentos desenvolvidos em laboratório didático; realização de relatórios para cada experimento.'
$ws.Cells.Item(18, 3).Value = 'Experimentos desenvolvidos em laboratório didático; realização de relatórios para cada experimento.'
$ws.Rows.Item(18).RowHeight = 60

# Row 19
$ws.Cells.Item(19, 1).Value = 'Norma de recuperação:'
$ws.Cells.Item(19, 2).Value = 'Média aritmética das notas obtidas nos relatórios. Será aprovado o aluno que obtiver nota final maior ou igual a 5,0.'
$ws.Cells.Item(19, 3).Value = 'Média aritmética das notas obtidas nos relatórios. Será aprovado o aluno que obtiver nota final maior ou igual a 5,0.'
$ws.Rows.Item(19).RowHeight = 60

# Row 20
$ws.Cells.Item(20, 1).Value = 'Bibliografia:'
$ws.Cells.Item(20, 2).Value = 'Devido às características práticas da disciplina, não será oferecida recuperação.'
$ws.Cells.Item(20, 3).Value = 'Devido às características práticas da disciplina, não será oferecida recuperação.'
$ws.Rows.Item(20).RowHeight = 120

# Row 21
$ws.Cells.Item(21, 1).Value = 'Requisitos:'
$ws.Cells.Item(21, 2).ClearContents() | Out-Null
$ws.Cells.Item(21, 3).ClearContents() | Out-Null

# Row 22
$ws.Cells.Item(22, 1).ClearContents() | Out-Null
$ws.Cells.Item(22, 2).Value = 'LOM3081 -  Introdução à Mecânica dos Sólidos  (Requisito fraco)
'
$ws.Cells.Item(22, 3).Value = 'LOM3081 -  Introdução à Mecânica dos Sólidos  (Requisito fraco)
'
$ws.Rows.Item(22).RowHeight = 30
